# Removed some not needed components.
# The BOM had two resistor board-reference lists that each still listed a
# resistor which is no longer populated on the board: R58 (1k resistors,
# row 31) and R56 (100k resistors, row 36). Drop them from the comma
# separated reference lists; all the quantities/pricing columns in the
# sheet are formulas driven off these two cells, so they recalc on their
# own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component List")

# --- Row 31 (1k resistors): remove "R58" -------------------------------
# C31 uses per-character (rich text) colouring for a few of the board
# references (R39/R59 green, R64 red) - set the plain value first and then
# re-apply the original colours to the surviving runs so the formatting
# for the untouched references is preserved, matching how Excel keeps
# the unaffected run formatting when you edit only part of a cell's text.
$c31 = $ws.Range("C31")
$c31.Value = "R10,R13,R16,R21,R23,R24,R29,R39,R50,R51,R57,R59,R62,R64"

$green = 5287936   # RGB(0,176,80)
$black = 0         # RGB(0,0,0)
$red = 255         # RGB(255,0,0)

$c31.Characters(29, 3).Font.Color = $green    # R39
$c31.Characters(32, 13).Font.Color = $black   # ,R50,R51,R57,
$c31.Characters(45, 3).Font.Color = $green    # R59
$c31.Characters(48, 5).Font.Color = $black    # ,R62,
$c31.Characters(53, 3).Font.Color = $red      # R64

# --- Row 36 (100k resistors): remove "R56" ------------------------------
$ws.Range("C36").Value = "R11,R14,R17,R35,R37,R38,R48,R49,R55"
